# Refactor generatePackage.js to improve metadata handling and enhance API mapping logic
# Update the header labels in the "Hoja1" sheet to use clearer column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename header cells: B1 "Meta" -> "MetadataType", D1 "Api" -> "ApiName"
$ws.Range("B1").Value = "MetadataType"
$ws.Range("D1").Value = "ApiName"

# Update the active selection on the sheet to match the saved view state
$ws.Activate()
$ws.Range("C15").Select()
